# The deck's slide/notes master theme ("Integral") is swapped back to the
# default "Office Theme" palette - dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink,
# in MsoThemeColorSchemeIndex order (1..12). PowerPoint's RGB() stores
# colors as 0x00BBGGRR, i.e. R + G*256 + B*65536.
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = 0x000000   # dk1     000000
$colorScheme.Colors(2).RGB  = 0xFFFFFF   # lt1     FFFFFF
$colorScheme.Colors(3).RGB  = 0x6A5444   # dk2     44546A
$colorScheme.Colors(4).RGB  = 0xE6E6E7   # lt2     E7E6E6
$colorScheme.Colors(5).RGB  = 0xD59B5B   # accent1 5B9BD5
$colorScheme.Colors(6).RGB  = 0x317DED   # accent2 ED7D31
$colorScheme.Colors(7).RGB  = 0xA5A5A5   # accent3 A5A5A5
$colorScheme.Colors(8).RGB  = 0x00C0FF   # accent4 FFC000
$colorScheme.Colors(9).RGB  = 0xC47244   # accent5 4472C4
$colorScheme.Colors(10).RGB = 0x47AD70   # accent6 70AD47
$colorScheme.Colors(11).RGB = 0xC16305   # hlink   0563C1
$colorScheme.Colors(12).RGB = 0x724F95   # folHlink 954F72
